$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset for the "falling" sheet was re-spliced: three fresh samples were
# inserted at the front (rows 2-4), the previously-existing samples shifted
# down by three rows (old rows 2-18 -> new rows 5-21), and the former last
# row (22) dropped so the sheet keeps the same number of data rows (20).

# Step 1: shift the existing sensor readings (columns C:H) for rows 2-18 down
# to rows 5-21, preserving the timestamp/label columns (A:B) as-is.
$src = $ws.Range("C2:H18")
$dst = $ws.Range("C5")
$src.Copy($dst)

# Step 2: write the three brand-new samples into rows 2-4 (columns C:H).
$ws.Cells.Item(2,3).Value = -0.4886150360107422
$ws.Cells.Item(2,4).Value = 1.498652458190918
$ws.Cells.Item(2,5).Value = -0.1321379840373993
$ws.Cells.Item(2,6).Value = -0.2070114476715818
$ws.Cells.Item(2,7).Value = -0.2780065764399136
$ws.Cells.Item(2,8).Value = 0.06705144135391006

$ws.Cells.Item(3,3).Value = -0.7675657272338867
$ws.Cells.Item(3,4).Value = 1.561143398284912
$ws.Cells.Item(3,5).Value = -0.3004561066627502
$ws.Cells.Item(3,6).Value = -0.1988007093177122
$ws.Cells.Item(3,7).Value = -0.2540031636462492
$ws.Cells.Item(3,8).Value = 0.1641969842945828

$ws.Cells.Item(4,3).Value = -0.6989822387695312
$ws.Cells.Item(4,4).Value = 1.441655874252319
$ws.Cells.Item(4,5).Value = -0.3177179098129272
$ws.Cells.Item(4,6).Value = -0.1353515688987338
$ws.Cells.Item(4,7).Value = -0.4769509890500239
$ws.Cells.Item(4,8).Value = -0.3330293473075421

# Step 3: the data that used to live in row 22 is now redundant (it was
# shifted into row 21's old slot by the splice), so remove the trailing row
# entirely. This also shrinks the sheet dimension from H22 to H21.
$ws.Rows.Item(22).Delete()
